$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "kerolsgamel2000@gmail.com"
$ws.Range("C2").Value = "https://github.com/kerolsgamel/Team-3-OSS"
$ws.Range("A2").Value = "كرلس جميل سامي بلامون"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:kerolsgamel2000@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/kerolsgamel/Team-3-OSS")

[void]$ws.Range("E6").Select()
